# "add Feedback form and Game Rules on first page" -- the concrete OOXML change this maps to
# is: 11 new English word-pairs appended to the "english" sheet (ending up as rows 13-23),
# renumbering of the existing English IDs (1000-1011 -> 1001-1012), widened B/C columns,
# extended conditional formatting, and the "english" tab becoming the active/front sheet
# (with fresh selections on both sheets).

$wb = $excel.ActiveWorkbook
$wsChinese = $wb.Worksheets.Item("chinese")
$wsEnglish = $wb.Worksheets.Item("english")

# --- 1. Renumber existing English rows' IDs (A1:A12): 1000..1011 -> 1001..1012 ---
for ($r = 1; $r -le 12; $r++) {
    $cell = $wsEnglish.Range("A$r")
    $cell.Value = $cell.Value() + 1
}

# --- 2. Append 11 new word-pair rows, each inserted via Rows.Insert() so the new row
#        inherits the existing body style (s="1") from the row above it, exactly like
#        Excel does natively for a typed/appended row. ---

# toothbrush / toothpaste ends up directly above Trump / Biden (row 13 & 14)
$wsEnglish.Rows(13).Insert()
$wsEnglish.Range("A13").Value = 1014
$wsEnglish.Range("B13").Value = "Trump"
$wsEnglish.Range("C13").Value = "Biden"

$wsEnglish.Rows(13).Insert()
$wsEnglish.Range("A13").Value = 1013
$wsEnglish.Range("B13").Value = "toothbrush"
$wsEnglish.Range("C13").Value = "toothpaste"
$wsEnglish.Range("A14").Value = 1014

$wsEnglish.Rows(15).Insert()
$wsEnglish.Range("A15").Value = 1015
$wsEnglish.Range("B15").Value = "eyebrow"
$wsEnglish.Range("C15").Value = "eyelash"

$wsEnglish.Rows(16).Insert()
$wsEnglish.Range("A16").Value = 1016
$wsEnglish.Range("B16").Value = "Spiderman"
$wsEnglish.Range("C16").Value = "Batman"

$wsEnglish.Rows(17).Insert()
$wsEnglish.Range("A17").Value = 1017
$wsEnglish.Range("C17").Value = "Cinderella"
$wsEnglish.Range("B17").Value = "Ugly Duckling"

# bread / cake ends up directly above lip balm / lipstick (row 18 & 19)
$wsEnglish.Rows(18).Insert()
$wsEnglish.Range("A18").Value = 1019
$wsEnglish.Range("B18").Value = "lip balm"

$wsEnglish.Rows(18).Insert()
$wsEnglish.Range("A18").Value = 1018
$wsEnglish.Range("B18").Value = "bread"
$wsEnglish.Range("C18").Value = "cake"
$wsEnglish.Range("A19").Value = 1019
$wsEnglish.Range("C19").Value = "lipstick"

$wsEnglish.Rows(20).Insert()
$wsEnglish.Range("A20").Value = 1020
$wsEnglish.Range("B20").Value = "tissue paper"
$wsEnglish.Range("C20").Value = "handkerchief"

$wsEnglish.Rows(21).Insert()
$wsEnglish.Range("A21").Value = 1021
$wsEnglish.Range("B21").Value = "chilli"
$wsEnglish.Range("C21").Value = "wasabi"

$wsEnglish.Rows(22).Insert()
$wsEnglish.Range("A22").Value = 1022
$wsEnglish.Range("B22").Value = "security guard"
$wsEnglish.Range("C22").Value = "bodyguard"

$wsEnglish.Rows(23).Insert()
$wsEnglish.Range("A23").Value = 1023
$wsEnglish.Range("B23").Value = "mirror"
$wsEnglish.Range("C23").Value = "glass"

# --- 3. Widen columns B and C to fit the new (longer) words. ---
$wsEnglish.Columns("B").ColumnWidth = 13.833333333333334
$wsEnglish.Columns("C").ColumnWidth = 13.0

# --- 4. Extend the duplicate-values conditional formatting to cover the new rows; the
#        priorities bump from 1/2 to 3/4 (matching the rule-creation counter Excel keeps
#        per sheet), while the two rules stay merged under one sqref / keep their dxfIds. ---
$fcs = $wsEnglish.Range("B1:C12").FormatConditions
$fc1 = $fcs.Item(1)
$fc2 = $fcs.Item(2)
$fc1.ModifyAppliesToRange($wsEnglish.Range("B1:C23"))
$fc2.ModifyAppliesToRange($wsEnglish.Range("B1:C23"))
$fc1.Priority = 3
$fc2.Priority = 4

# --- 5. Update view state: "chinese" loses its old scroll position / selection and is no
#        longer the front tab; "english" becomes the active tab with a fresh selection. ---
$wsChinese.Activate()
$wsChinese.Range("D27").Select()

$wsEnglish.Activate()
$wsEnglish.Range("I28").Select()
